$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 and 40 swap places (EnergySwap moves to row 39, LidoDAOToken moves to row 40)
# along with their updated Price / Volume(1h) values
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '22.83'
$ws.Range("E39").Value = '  +11.75%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +2.21%  '

# Update Price (D) and Volume(1h) (E) columns for other rows with changed values
$ws.Range("D2").Value = '43.037.72'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '2.375.54'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '302.01'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '96.63'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").Value = '34.18'
$ws.Range("E10").Value = '  -1.42%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").Value = '18.28'
$ws.Range("E13").Value = '  -4.55%  '
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '2.747.42'
$ws.Range("E15").Value = '  +2.40%  '
$ws.Range("D16").Value = '2.377.53'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("D18").Value = '43.005.79'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '12.15'
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").Value = '0.0₃0887'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = '68.08'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '235.33'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("E24").Value = '  -2.23%  '
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '24.91'
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").Value = '31.43'
$ws.Range("E30").Value = '  -3.20%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").Value = '0.0738'
$ws.Range("E33").Value = '  +5.09%  '
$ws.Range("D34").Value = '17.45'
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").Value = '1.86'
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("E36").Value = '  +5.00%  '
$ws.Range("E37").Value = '  -2.99%  '
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").Value = '113.05'
$ws.Range("E42").Value = '  -32.10%  '
$ws.Range("D43").Value = '1.951.81'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '2.73'
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").Value = '9.14'
$ws.Range("E47").Value = '  -11.69%  '
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").Value = '52.31'
$ws.Range("E49").Value = '  -2.23%  '
$ws.Range("D50").Value = '72.31'
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("E51").Value = '  +0.61%  '
